# The deck's slide master / presentation theme (ppt/theme/theme2.xml) is
# currently the "Integral" theme. The author's edit swaps the colour
# palettes that live in ppt/theme/theme1.xml and ppt/theme/theme2.xml, so
# the master's theme ends up using the stock "Office Theme" colours
# (what used to be theme1.xml) while the old "Integral" colours move to
# theme1.xml (used only by the notes master).
#
# PowerPoint's ColorScheme object is the supported automation surface for
# recolouring a theme, so recreate the swap by pushing the "Office Theme"
# RGB values onto the (only reachable / in-use) master colour scheme.
# dk1/lt1 (black/white) are identical in both palettes, so only the other
# ten slots need updating.

function ConvertTo-VbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Item(3).RGB  = ConvertTo-VbaRgb "44546A"   # dk2
$cs.Item(4).RGB  = ConvertTo-VbaRgb "E7E6E6"   # lt2
$cs.Item(5).RGB  = ConvertTo-VbaRgb "5B9BD5"   # accent1
$cs.Item(6).RGB  = ConvertTo-VbaRgb "ED7D31"   # accent2
$cs.Item(7).RGB  = ConvertTo-VbaRgb "A5A5A5"   # accent3
$cs.Item(8).RGB  = ConvertTo-VbaRgb "FFC000"   # accent4
$cs.Item(9).RGB  = ConvertTo-VbaRgb "4472C4"   # accent5
$cs.Item(10).RGB = ConvertTo-VbaRgb "70AD47"   # accent6
$cs.Item(11).RGB = ConvertTo-VbaRgb "0563C1"   # hyperlink
$cs.Item(12).RGB = ConvertTo-VbaRgb "954F72"   # followed hyperlink
